$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update title (D6) and link (E6)
$ws.Range("D6").Value = "Object Detection이란? Object Detection 용어정리"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Object-Detection%EC%9D%B4%EB%9E%80-Object-Detection-%EC%9A%A9%EC%96%B4%EC%A0%95%EB%A6%AC"

# Row 26: update title (D26) only
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 46: update title (D46) and link (E46)
$ws.Range("D46").Value = "간염 (Hepatitis)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/450"
